$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.82231633333333
$ws.Range("H2").Value = 44.466949
$ws.Range("I2").Value = 0.0966878806285147
$ws.Range("J2").Value = 0.0966878806285147
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 40.37601439586056
$ws.Range("R2").Value = 363.384129562745
$ws.Range("S2").Value = 0.00448196716506411
$ws.Range("T2").Value = 0.00448196716506411
$ws.Range("G3").Value = 14.82231633333333
$ws.Range("H3").Value = 44.466949
$ws.Range("I3").Value = 0.0966878806285147
$ws.Range("J3").Value = 0.0966878806285147
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 602.0918374008343
$ws.Range("R3").Value = 5418.826536607508
$ws.Range("S3").Value = 0.0668356172832235
$ws.Range("T3").Value = 0.06683561728322349
$ws.Range("G4").Value = 14.82231633333333
$ws.Range("H4").Value = 44.466949
$ws.Range("I4").Value = 0.0966878806285147
$ws.Range("J4").Value = 0.0966878806285147
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 228.5495199038217
$ws.Range("R4").Value = 2056.945679134395
$ws.Range("S4").Value = 0.02537029618022709
$ws.Range("T4").Value = 0.02537029618022709
$ws.Range("I5").Value = 0.1609965995515919
$ws.Range("J5").Value = 0.1609965995515918
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 67.23077369080944
$ws.Range("R5").Value = 605.076963217285
$ws.Range("S5").Value = 0.00746299813571883
$ws.Range("T5").Value = 0.007462998135718829
$ws.Range("I6").Value = 0.1609965995515919
$ws.Range("J6").Value = 0.1609965995515918
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("S6").Value = 0.1112890989189519
$ws.Range("T6").Value = 0.1112890989189519
$ws.Range("I7").Value = 0.1609965995515919
$ws.Range("J7").Value = 0.1609965995515918
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 380.5616101467483
$ws.Range("R7").Value = 3425.054491320735
$ws.Range("S7").Value = 0.04224450249692113
$ws.Range("T7").Value = 0.04224450249692113
$ws.Range("G8").Value = 113.7974623333333
$ws.Range("H8").Value = 341.392387
$ws.Range("I8").Value = 0.7423155198198935
$ws.Range("J8").Value = 0.7423155198198935
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 309.9844770584372
$ws.Range("R8").Value = 2789.860293525935
$ws.Range("S8").Value = 0.034410039441583
$ws.Range("T8").Value = 0.034410039441583
$ws.Range("G9").Value = 113.7974623333333
$ws.Range("H9").Value = 341.392387
$ws.Range("I9").Value = 0.7423155198198935
$ws.Range("J9").Value = 0.7423155198198935
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 4622.524688246245
$ws.Range("R9").Value = 41602.72219421621
$ws.Range("S9").Value = 0.5131265228234598
$ws.Range("T9").Value = 0.5131265228234598
$ws.Range("G10").Value = 113.7974623333333
$ws.Range("H10").Value = 341.392387
$ws.Range("I10").Value = 0.7423155198198935
$ws.Range("J10").Value = 0.7423155198198935
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 1754.675503994432
$ws.Range("R10").Value = 15792.07953594988
$ws.Range("S10").Value = 0.1947789575548506
$ws.Range("T10").Value = 0.1947789575548506
